# Saldo_guide.xlsx update: refresh source extract from 2024-05-31 to 2024-06-04.
#  - rename the sheet to reflect the new extraction timestamp
#  - bump every "Dt. Referencia" (column G) date by one day (45446 -> 45447)
#  - for accounts whose "Vl. Projetado" (column E) had a pending/projected
#    amount, roll it into "Saldo Previsto" (column D) now that it has
#    settled, leaving E at 0 (H = D + E stays the same)
#  - a few accounts also got a refreshed "Saldo Previsto"/"Vl. Total" figure
#    from the new extract even though E was already 0

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet to match the new extraction run.
$ws.Name = "IClientBalance-20240604-082626-"

# 2) Every data row (2..257): advance the reference date by one day.
$lastRow = 257
for ($r = 2; $r -le $lastRow; $r++) {
    $g = $ws.Cells.Item($r, 7).Value2
    $ws.Cells.Item($r, 7).Value2 = $g + 1
}

# 3) Rows where the projected value (E) settled into the forecast balance
#    (D): new D = old D + old E (the exact figure, taken from the new
#    extract so it matches to the cent), E reset to 0. H (Vl. Total) was
#    already equal to D + E and is NOT re-derived here, so no floating
#    point drift is introduced by re-adding D and E.
$newD = @{
    5   = 11711.55
    8   = 4755.47
    15  = 14679.76
    17  = 5260.14
    42  = 6989.84
    57  = 2184.98
    59  = 9450.27
    98  = 9142.2000000000007
    103 = 24792.95
    107 = 29054.59
    131 = 3949.98
    141 = 32625.55
    155 = 855.12
    168 = 3335.98
    226 = 8179.5
    240 = 8472.92
}
foreach ($r in $newD.Keys) {
    $ws.Cells.Item($r, 4).Value2 = $newD[$r]
    $ws.Cells.Item($r, 5).Value2 = 0
}

# 4) Rows with a straight refreshed balance from the new extract (E was
#    already 0, D/H just get a new figure -- H mirrors the new D exactly,
#    same as it mirrored the old D before).
$ws.Cells.Item(51, 4).Value2 = 124.61
$ws.Cells.Item(51, 8).Value2 = 124.61

$ws.Cells.Item(111, 4).Value2 = 130968.58
$ws.Cells.Item(111, 8).Value2 = 130968.58

$ws.Cells.Item(136, 4).Value2 = 557.23
$ws.Cells.Item(136, 8).Value2 = 557.23

# 5) Collapse the selection back down to the top-left cell (matches the
#    saved workbook no longer keeping the whole used range selected).
$ws.Range("A1").Select() | Out-Null
